$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

$ws.Range("A19").Value = "commissioner2"
$ws.Range("B19").Value = "ADMINISTRATION"
$ws.Range("C19").Value = "Commissioner"
$ws.Range("D19").Value = "Ravindra Babu ~ ADM_Commissioner_2"
$ws.Range("E19").Value = "Forward to commissioner"

$ws.Range("D24").Select()
